$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: Adam Smith
$ws.Range("B5").Value = "Adam"
$ws.Range("C5").Value = "Smith"

# Move the active selection down to B6 (as it would after typing into B5:C5)
$ws.Range("B6").Select()

# Resize the workbook window (visual improvement mentioned in commit message)
$excel.ActiveWindow.Width = 20730
$excel.ActiveWindow.Height = 11760
